$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.572.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.23%  '
$ws.Range("D3").Value = '''2.557.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.81%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''509.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.47%  '
$ws.Range("D6").Value = '''146.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.71%  '
$ws.Range("D8").Value = '''0.569'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.95%  '
$ws.Range("D9").Value = '''2.571.71'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.91%  '
$ws.Range("D10").Value = '''6.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.74%  '
$ws.Range("E11").Value = '  -6.08%  '
$ws.Range("E12").Value = '  -4.84%  '
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("D14").Value = '''3.005.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.73%  '
$ws.Range("D15").Value = '''58.515.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.28%  '
$ws.Range("D16").Value = '''21.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.73%  '
$ws.Range("D17").Value = '''0.0000137'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.16%  '
$ws.Range("D18").Value = '''2.563.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").Value = '''347.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.23%  '
$ws.Range("E21").Value = '  -4.51%  '
$ws.Range("E22").Value = '  -4.46%  '
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").Value = '''60.60'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("E25").Value = '  -4.37%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '''2.670.90'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.58%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.160'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.29%  '
$ws.Range("D29").Value = '''0.0₃0802'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.15%  '
$ws.Range("E30").Value = '  -5.04%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").Value = '''5.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("E33").Value = '  -4.93%  '
$ws.Range("D34").Value = '''149.39'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("E35").Value = '  -5.89%  '
$ws.Range("E36").Value = '  -4.72%  '
$ws.Range("D37").Value = '''0.899'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.20%  '
$ws.Range("E38").Value = '  -6.03%  '
$ws.Range("D39").Value = '''0.844'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.71%  '
$ws.Range("D40").Value = '''36.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.12%  '
$ws.Range("E41").Value = '  -6.50%  '
$ws.Range("D42").Value = '''285.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.66%  '
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("D45").Value = '''0.997'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '''0.606'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.22%  '
$ws.Range("E47").Value = '  -5.07%  '
$ws.Range("D48").Value = '''19.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.61%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = '''10.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '''0.0229'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.50%  '
$ws.Range("D51").Value = '''4.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.13%  '
